$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (rename storageColorOption -> storageColorOptions)
$ws.Range("G1").Value = "label_storageColorOptions_for"
$ws.Range("H1").Value = "label_storageColorOptions_for_1"
$ws.Range("I1").Value = "label_storageColorOptions_internalText"
$ws.Range("J1").Value = "label_storageColorOptions_internalText_1"

# Update data row value (value is textual "2", not numeric -- keep it text
# with a leading apostrophe so Excel stores it as a string rather than a number)
$ws.Range("E2").Value = "'2"

# Adjust column widths for G, H, I, J (bump up by 1 character each)
# Note: COM ColumnWidth differs from raw OOXML stored width by a constant
# offset (~0.8333 for this workbook's default font), so we compensate here
# so the saved XML "width" attribute ends up exactly 31 / 33 / 40 / 42.
$ws.Range("G1").ColumnWidth = 30.166666666666668
$ws.Range("H1").ColumnWidth = 32.166666666666664
$ws.Range("I1").ColumnWidth = 39.166666666666664
$ws.Range("J1").ColumnWidth = 41.166666666666664
